# "fixed pyomo not running" — correct the Sheet1 boolean/override cells that were
# feeding bad values (9999 / 0.9) into the pyomo export, and point the selection
# at A6 instead of the old J7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A5 was TRUE, needs to be FALSE
$ws.Range("A5").Value = $false

# J5 held an IFS() array formula (style kept) — replace with a plain literal 1
$ws.Range("J5").Value = 1

# J6 / K6 were stray overrides (9999 / 0.9) — reset them to 1.
# Everything below (J7:J29) is driven off J6 via `=J{n-1}` shared formulas, so
# fixing J6 ripples the cached value down the whole chain automatically.
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 1

# Leave the selection on A6, matching where the fix was made.
$ws.Range("A6").Select() | Out-Null
